$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.211.24"
$ws.Range("E2").Value = "  -1.61%  "
$ws.Range("D3").Value = "2.249.47"
$ws.Range("E3").Value = "  -1.38%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.12"
$ws.Range("E5").Value = "  -1.14%  "
$ws.Range("E6").Value = "  -2.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.51"
$ws.Range("E7").Value = "  -6.21%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.618"
$ws.Range("E9").Value = "  -4.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.17"
$ws.Range("E10").Value = "  +2.62%  "
$ws.Range("E11").Value = "  -2.84%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.20"
$ws.Range("E12").Value = "  -2.24%  "
$ws.Range("E13").Value = "  -1.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.52"
$ws.Range("E14").Value = "  -3.70%  "
$ws.Range("E15").Value = "  -1.77%  "
$ws.Range("D16").Value = "2.254.13"
$ws.Range("E16").Value = "  -1.30%  "
$ws.Range("D17").Value = "42.134.39"
$ws.Range("E17").Value = "  -1.59%  "
$ws.Range("D18").Value = "0.0₃0991"
$ws.Range("E18").Value = "  -0.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.45"
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.13"
$ws.Range("E21").Value = "  +3.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "231.64"
$ws.Range("E22").Value = "  -1.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.88"
$ws.Range("E23").Value = "  +36.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.41"
$ws.Range("E25").Value = "  +0.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.62"
$ws.Range("E26").Value = "  -4.85%  "
$ws.Range("E27").Value = "  -1.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.24"
$ws.Range("E28").Value = "  +2.84%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "169.02"
$ws.Range("E29").Value = "  +0.46%  "
$ws.Range("E30").Value = "  -1.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0827"
$ws.Range("E31").Value = "  -3.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.123"
$ws.Range("E32").Value = "  -0.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.10"
$ws.Range("E33").Value = "  +3.40%  "
$ws.Range("E34").Value = "  -1.89%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.22"
$ws.Range("E35").Value = "  +9.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.48"
$ws.Range("E36").Value = "  -1.93%  "
$ws.Range("E37").Value = "  +3.62%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "13.89"
$ws.Range("E38").Value = "  +1.41%  "
$ws.Range("E39").Value = "  -3.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.78"
$ws.Range("E40").Value = "  -1.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "62.71"
$ws.Range("E41").Value = "  +2.16%  "
$ws.Range("E42").Value = "  -2.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "106.83"
$ws.Range("E43").Value = "  -4.56%  "
$ws.Range("E44").Value = "  +1.38%  "
$ws.Range("E45").Value = "  -2.53%  "
$ws.Range("E46").Value = "  -0.37%  "
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.12"
$ws.Range("E47").Value = "  -2.99%  "
$ws.Range("B48").Value = "TrustWalletToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.17"
$ws.Range("E48").Value = "  +0.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.30"
$ws.Range("E49").Value = "  +2.35%  "
$ws.Range("E50").Value = "  -10.24%  "
$ws.Range("B51").Value = "SynthetixNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.10"
$ws.Range("E51").Value = "  -3.41%  "
